$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The very first paragraph in the document is empty (just a paragraph
#    mark). Its paragraph-mark run properties need a `w:hint="cs"` added to
#    the `<w:rFonts .../>` element. The Word object model has no dedicated
#    "hint" property, so we target that exact paragraph and rewrite its
#    OOXML via InsertXML, keeping every other property identical.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00A93B23" w:rsidRDefault="00A93B23"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="cs"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl/></w:rPr></w:pPr></w:p>
'@
$firstPara.Range.InsertXML($firstParaXml)

# ---------------------------------------------------------------------------
# 2) "...קרא לטבלה file." -> "...קרא לטבלה Physician." (single occurrence).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("file", $true, $true, $false, $false, $false, $true, 1, $false, "Physician", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) The paragraph that ends "...נניח Desktop. " has a hidden `_GoBack`
#    bookmark sitting right after the "Desktop" / ". " runs. It needs to move
#    so the bookmark wraps *before* those two runs instead of after them
#    (same runs, same formatting - only their position relative to the
#    bookmark changes). `_GoBack` isn't reachable through Document.Bookmarks
#    in this host, so the paragraph is located by its text and rewritten
#    in place with InsertXML, reusing the exact original run XML.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Desktop*") {
        $targetPara = $candidate
        break
    }
}

$desktopParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007B164A" w:rsidRDefault="005E6466"><w:pPr><w:pStyle w:val="NormalWeb"/><w:bidi/><w:spacing w:before="0" w:after="0"/><w:ind w:firstLine="360"/><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl/></w:rPr><w:t>במקום ה-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>path</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl/></w:rPr><w:t xml:space="preserve">, רשום את המיקום, נניח </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>Desktop</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>
'@

if ($targetPara -ne $null) {
    $targetPara.Range.InsertXML($desktopParaXml)
}

Write-Output "done"
